# Update gh-pages to output generated at 456a3b4
# Bumps the "想去人数" (want-to-go count) numbers in the F column
# on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1243
$ws1.Range("F4").Value = 2734
$ws1.Range("F5").Value = 244

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1243
$ws4.Range("F6").Value = 2734
$ws4.Range("F8").Value = 244
